$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 21: relabel "Stim.ESS2.MCB" -> "Stim.ESS2.VFmode" and invert the 0/1 pattern ---
$ws.Range("A21").Value = "Stim.ESS2.VFmode"
$ws.Range("C21:J21").Value = 1
$ws.Range("K21:AL21").Value = 0

# --- Row 28: battery-trip formula now keys off VFmode flag (C21=1 instead of C21=0) ---
$ws.Range("D28:AL28").Formula = "=IF(OR(D16=0,D16=0),0,IF(D21=1,""GF"",D17*250))"
$ws.Range("C28").Formula = "=IF(OR(C16=0,C16=0),0,IF(C21=1,""GF"",C17*250))"

# --- Row 35: same condition update on the Q-side battery-trip row ---
$ws.Range("C35").Formula = "=IF(OR(C16=0,C16=0),0,IF(C21=1,""GF"",C18*250))"
$ws.Range("D35:AL35").Formula = "=IF(OR(D16=0,D16=0),0,IF(D21=1,""GF"",D18*250))"

# --- Column A autosized for the longer label, selection left on the row 35 block ---
$ws.Columns("A").ColumnWidth = 20.6
$ws.Activate()
$ws.Range("C35:AL35").Select()

Write-Output "done"
